# Progression chapitre 5 Lidar
#
# 1) Add the new data row (date 09/08/18 -> 43321 in B5, value 120 in C5),
#    copying the date number-format from B4 so B5 gets the same style.
# 2) Convert the chart from a stacked column chart to a line chart
#    (with no markers, not smoothed).
# 3) Update the active cell selection on the sheet.
# 4) Nudge the workbook window position to match the saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New data row -------------------------------------------------
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B5").Value = 43321
$ws.Range("C5").Value = 120

# --- 2. Chart type: stacked bar -> line -------------------------------
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$chart.ChartType = 4   # xlLine

for ($i = 1; $i -le $chart.SeriesCollection().Count; $i++) {
    $s = $chart.SeriesCollection($i)
    $s.MarkerStyle = -4142   # xlMarkerStyleNone
    $s.Smooth = $false
}

# --- 3. Selection -------------------------------------------------------
$ws.Range("C15").Select() | Out-Null

# --- 4. Window position (best effort) -----------------------------------
$win = $wb.Windows.Item(1)
$win.Left = 7920
$win.Top = 2940
